# Update "想去人数" (want-to-go count, column F) figures across the
# 展览 / 演出 / 全部类型 sheets to match the refreshed scrape output.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 5972
$ws.Range("F5").Value = 5972
$ws.Range("F7").Value = 3021
$ws.Range("F8").Value = 1294
$ws.Range("F9").Value = 447
$ws.Range("F12").Value = 34
$ws.Range("F13").Value = 321
$ws.Range("F14").Value = 4447
$ws.Range("F15").Value = 4447
$ws.Range("F16").Value = 102
$ws.Range("F17").Value = 95
$ws.Range("F21").Value = 80
$ws.Range("F22").Value = 6890
$ws.Range("F23").Value = 6890
$ws.Range("F24").Value = 237
$ws.Range("F25").Value = 110
$ws.Range("F26").Value = 475
$ws.Range("F27").Value = 1271
$ws.Range("F28").Value = 6267
$ws.Range("F29").Value = 1648
$ws.Range("F31").Value = 1971
$ws.Range("F32").Value = 6033
$ws.Range("F36").Value = 88
$ws.Range("F37").Value = 436
$ws.Range("F38").Value = 6003
$ws.Range("F44").Value = 2419
$ws.Range("F47").Value = 1010
$ws.Range("F49").Value = 357
$ws.Range("F50").Value = 2084
$ws.Range("F52").Value = 1034

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 208
$ws.Range("F7").Value = 31

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 5972
$ws.Range("F5").Value = 5972
$ws.Range("F7").Value = 3021
$ws.Range("F8").Value = 1294
$ws.Range("F9").Value = 447
$ws.Range("F12").Value = 208
$ws.Range("F13").Value = 321
$ws.Range("F14").Value = 4447
$ws.Range("F15").Value = 4447
$ws.Range("F16").Value = 102
$ws.Range("F17").Value = 95
$ws.Range("F21").Value = 80
$ws.Range("F22").Value = 6890
$ws.Range("F23").Value = 6890
$ws.Range("F24").Value = 237
$ws.Range("F25").Value = 110
$ws.Range("F26").Value = 475
$ws.Range("F27").Value = 1271
$ws.Range("F28").Value = 6267
$ws.Range("F29").Value = 1648
$ws.Range("F30").Value = 31
$ws.Range("F32").Value = 1971
$ws.Range("F33").Value = 6033
$ws.Range("F37").Value = 88
$ws.Range("F38").Value = 436
$ws.Range("F39").Value = 6003
$ws.Range("F46").Value = 2419
$ws.Range("F47").Value = 1010
$ws.Range("F48").Value = 357
$ws.Range("F49").Value = 2084
$ws.Range("F51").Value = 1034
